$d = $word.ActiveDocument

# Paragraphs that need a first-line indent of 360 twips (18 points) added
# to their paragraph properties (w:pPr/w:ind w:firstLine="360").
$targets = @(
    "The first week of development are more focused",
    "Project management side, ",
    "The future dedicated branches are for large-based",
    "For the boilerplate template, I successfully create",
    "My plan for the first agile iteration"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    foreach ($prefix in $targets) {
        if ($text.StartsWith($prefix)) {
            $p.Range.ParagraphFormat.FirstLineIndent = 18
            break
        }
    }
}
